$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: the "Objetivos:" value cells (already exist) get the real objectives text
$ws.Range("B10").Value = "Introduzir ao aluno a teoria de propriedades elétricas, térmicas, magnéticas e óticas de materiais sólidos, levando emconta o aspecto microscópico da estrutura do material. Dá-se ênfase à aplicação do material de acordo com aspropriedades que ele apresenta."
$ws.Range("C10").Value = "Introduzir ao aluno a teoria de propriedades elétricas, térmicas, magnéticas e óticas de materiais sólidos, levando emconta o aspecto microscópico da estrutura do material. Dá-se ênfase à aplicação do material de acordo com aspropriedades que ele apresenta."

# Row 13: drop the stray label in A13; B13/C13 keep their style and get the first
# professor entry
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C13").Value = "5840726 - Cristina Bormio Nunes"

# Row 14: drop the stray label in A14; create B14/C14 (style copied from B13/C13)
# with the second professor entry
$ws.Range("A14").Clear()
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B14").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C14").Value = "1341653 - Maria José Ramos Sandim"

# Row 15 ("Programa resumido:") now carries the short-syllabus paragraph
$ws.Range("B15").Value = "PROPRIEDADES ELETRÔNICAS: Condutividade elétrica em metais, semicondutores e isolantes.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão, Condutividade e Tensões Térmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Lasers. Aplicações."
$ws.Range("C15").Value = "PROPRIEDADES ELETRÔNICAS: Condutividade elétrica em metais, semicondutores e isolantes.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão, Condutividade e Tensões Térmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Lasers. Aplicações."

# Row 17 ("Programa:") did not have B/C cells before; copy style from B15/C15 and
# fill in the full syllabus paragraph
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B17").Value = "PROPRIEDADES ELETRÔNICAS:Teoria do Elétron Livre em Metais. Níveis de Energia em Sólidos. Condutividade.Supercondutividade. Semicondutividade. Isolantes (Dielétricos). Aplicações.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão Térmica. Condutividade Térmica. TensõesTérmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Fotocondutividade. Luminescência. Lasers.Fibra Ótica. Danos por Radiação. Aplicações."
$ws.Range("C17").Value = "PROPRIEDADES ELETRÔNICAS:Teoria do Elétron Livre em Metais. Níveis de Energia em Sólidos. Condutividade.Supercondutividade. Semicondutividade. Isolantes (Dielétricos). Aplicações.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão Térmica. Condutividade Térmica. TensõesTérmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Fotocondutividade. Luminescência. Lasers.Fibra Ótica. Danos por Radiação. Aplicações."

# Row 18 becomes a lone "Syllabus:" label row — drop its old B/C values
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()

# Row 19 becomes a lone "Avaliação:" label row — drop its old B/C values
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# Row 20 ("Método:") now carries the exam-method sentence
$ws.Range("B20").Value = "Aplicação de duas provas em sala de aula, P1 e P2."
$ws.Range("C20").Value = "Aplicação de duas provas em sala de aula, P1 e P2."

# Row 21 ("Critério:") now carries the NF formula text
$ws.Range("B21").Value = "A nota final (NF) antes da recuperação será
NF = (P1 + 2*P2)/3"
$ws.Range("C21").Value = "A nota final (NF) antes da recuperação será
NF = (P1 + 2*P2)/3"

# Row 22 ("Norma de recuperação:") did not have B/C cells before; copy style from
# B21/C21 and fill in the recovery-rule text
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B22").Value = "A recuperação final consta da aplicação de uma prova PR A média final (MF) após a recuperação será
calculada como: MF = (NF + PR)/2"
$ws.Range("C22").Value = "A recuperação final consta da aplicação de uma prova PR A média final (MF) após a recuperação será
calculada como: MF = (NF + PR)/2"

# Row 23 becomes the "Bibliografia:" label row; add A23 label (style copied from A22)
# and set the bibliography text in B23/C23 (cells already existed)
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "[1] Callister, W. D., Materials Science and Engineering, An Introduction, John Wiley &Sons, 1990.
[2] Shakelford, J. F., Introduction to Materials Science for Engineers,,Prentice Hall, 4a. edição, 1996.
[3] Jastrzebski, Z. D., The Nature and Properties of Engineering Materials, John Wiley & Sons, 3a. edição, 1987.
[4] Solymar, L. and Walsh, D., Lectures On the Electrical Properties of Materials, Oxford Science Publications, 5a.
edição, 1993.
[5] Kittel, C., Introduction to Solid State Physics, John Wiley, 7a. edição, 1996"
$ws.Range("C23").Value = "[1] Callister, W. D., Materials Science and Engineering, An Introduction, John Wiley &Sons, 1990.
[2] Shakelford, J. F., Introduction to Materials Science for Engineers,,Prentice Hall, 4a. edição, 1996.
[3] Jastrzebski, Z. D., The Nature and Properties of Engineering Materials, John Wiley & Sons, 3a. edição, 1987.
[4] Solymar, L. and Walsh, D., Lectures On the Electrical Properties of Materials, Oxford Science Publications, 5a.
edição, 1993.
[5] Kittel, C., Introduction to Solid State Physics, John Wiley, 7a. edição, 1996"

# Row 24 becomes a lone "Requisitos:" label row; add A24 (style copied from A22) and
# drop the old B24/C24 values
$ws.Range("A22").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A24").Value = "Requisitos:"
$ws.Range("B24").Clear()
$ws.Range("C24").Clear()

# Rows 25/26: brand-new rows holding the two weak-requirement lines (style copied
# from B23/C23)
$ws.Range("B23").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("B23").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B25").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)
"
$ws.Range("C25").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)
"
$ws.Range("B26").Value = "LOM3109 -  Fundamentos da Física Moderna  (Requisito fraco)
"
$ws.Range("C26").Value = "LOM3109 -  Fundamentos da Física Moderna  (Requisito fraco)
"

# --- Row heights ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30

# --- Column layout tidy-up: column A should only span col 1, not 1:2 ---
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
